# Automated map update (2025-07-21 07:30:56)
# Appends the new incident row (row 25) to the INCO sheet, matching the
# existing sheet's column layout:
#   A Caso | B F. De Reclamo | C Direccion | D Comuna | E OT
#   F Proveedor Asignado | G Estado | H Observaciones | I Attachments
#   J Tipo de tarea | K Equipo | L Tipo de Elemento
#   M Coordenada_X | N Coordenada_Y | O Operacion | P Zona

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

# Columns A, B, D and E hold numeric-/date-looking text in every existing
# row (e.g. "-20", "2/6/2024", "1", "780027603"), so force them to the
# Text number format first - otherwise Excel would auto-coerce these
# values into a real number or date serial instead of keeping literal text.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("E$row").NumberFormat = "@"

$ws.Range("A$row").Value = "-522"
$ws.Range("B$row").Value = "7/21/2025"
$ws.Range("C$row").Value = "Uruguay 1090"
$ws.Range("D$row").Value = "2"
$ws.Range("E$row").Value = "808430941"
$ws.Range("F$row").Value = "INCO"
$ws.Range("G$row").Value = "Pendiente"
$ws.Range("H$row").Value = "Reclaman columna corroida y rienda fuera de norma pero no se ve en la foto."
$ws.Range("I$row").Value = 1
$ws.Range("J$row").Value = "Cambio"
$ws.Range("K$row").Value = "Sin equipos"
$ws.Range("L$row").Value = "Terminal"
$ws.Range("M$row").Value = -58.387175
$ws.Range("N$row").Value = -34.596
$ws.Range("O$row").Value = "Recoleta"
$ws.Range("P$row").Value = "Capital Sur"
